$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to insert before the current row 2 (shifts existing data down)
$topRows = @(
    @(-2.236848592758179, 3.817810773849488, 1.356045484542847),
    @(-2.285426902770996, 3.86536750793457, 1.423394083976746),
    @(-2.325414371490478, 3.847443521022797, 1.428170895576477),
    @(-2.198427677154541, 3.85212025642395, 1.377349805831909),
    @(-2.239168739318848, 3.77459921836853, 1.236697590351105),
    @(-2.284214735031127, 3.738507509231567, 1.190666794776917),
    @(-2.274757814407349, 3.780967509746552, 1.303351855278015),
    @(-2.355447578430176, 3.771934032440186, 1.401894807815552)
)

# New rows to append after the existing data
$bottomRows = @(
    @(1.128712320327757, 4.191518974304199, 1.193370014429091),
    @(0.9585402488708525, 4.21038007736206, 1.021172881126405)
)

# Insert rows at row 2, pushing existing data down
$insertRange = $ws.Range("A2:C" + (1 + $topRows.Count))
$insertRange.Insert(-4121) # xlShiftDown

for ($i = 0; $i -lt $topRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $topRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $topRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $topRows[$i][2]
}

# Clear any formatting the insert may have inherited from the header row so
# the new rows stay unstyled, matching the rest of the data rows
$newRowsRange = $ws.Range("A2:C" + (1 + $topRows.Count))
$newRowsRange.ClearFormats()

# Determine the first empty row after existing data and append bottom rows
$lastRow = 1 + $topRows.Count + 20
for ($i = 0; $i -lt $bottomRows.Count; $i++) {
    $r = $lastRow + 1 + $i
    $ws.Cells.Item($r, 1).Value = $bottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $bottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $bottomRows[$i][2]
}
